$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy formatting from column F (the old column D, now shifted) into the two new columns D:E
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)

# Update data: two new quarter columns (D, E) plus refreshed figures in the existing columns
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(7, 7).Value = 43190
$ws.Cells.Item(7, 8).Value = 43100
$ws.Cells.Item(7, 9).Value = 43008
$ws.Cells.Item(7, 10).Value = 42916
$ws.Cells.Item(7, 11).Value = 42825
$ws.Cells.Item(7, 12).Value = 42735
$ws.Cells.Item(7, 13).Value = 42643
$ws.Cells.Item(8, 4).Value = 1929500
$ws.Cells.Item(8, 5).Value = 1894100
$ws.Cells.Item(8, 6).Value = 1951200
$ws.Cells.Item(8, 7).Value = 1980500
$ws.Cells.Item(8, 8).Value = 1935000
$ws.Cells.Item(8, 9).Value = 1937600
$ws.Cells.Item(8, 10).Value = 1926700
$ws.Cells.Item(8, 11).Value = 1995700
$ws.Cells.Item(8, 12).Value = 2018000
$ws.Cells.Item(8, 13).Value = 1964600
$ws.Cells.Item(9, 4).Value = 1546100
$ws.Cells.Item(9, 5).Value = 1503500
$ws.Cells.Item(9, 6).Value = 1518400
$ws.Cells.Item(9, 7).Value = 1532000
$ws.Cells.Item(9, 8).Value = 1488500
$ws.Cells.Item(9, 9).Value = 1495800
$ws.Cells.Item(9, 10).Value = 1459200
$ws.Cells.Item(9, 11).Value = 1533500
$ws.Cells.Item(9, 12).Value = 1516600
$ws.Cells.Item(9, 13).Value = 1475800
$ws.Cells.Item(10, 4).Value = 383400
$ws.Cells.Item(10, 5).Value = 390600
$ws.Cells.Item(10, 6).Value = 432800
$ws.Cells.Item(10, 7).Value = 448500
$ws.Cells.Item(10, 8).Value = 446500
$ws.Cells.Item(10, 9).Value = 441800
$ws.Cells.Item(10, 10).Value = 467500
$ws.Cells.Item(10, 11).Value = 462200
$ws.Cells.Item(10, 12).Value = 501400
$ws.Cells.Item(10, 13).Value = 488800
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"
$ws.Cells.Item(12, 12).Value = "NA"
$ws.Cells.Item(12, 13).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(14, 4).Value = 203700
$ws.Cells.Item(14, 5).Value = -2700
$ws.Cells.Item(14, 6).Value = 67600
$ws.Cells.Item(14, 7).Value = 8500
$ws.Cells.Item(14, 8).Value = 7700
$ws.Cells.Item(14, 9).Value = 32800
$ws.Cells.Item(14, 10).Value = 5800
$ws.Cells.Item(14, 11).Value = 26300
$ws.Cells.Item(14, 12).Value = -300
$ws.Cells.Item(14, 13).Value = 9300
$ws.Cells.Item(15, 4).Value = 5200
$ws.Cells.Item(15, 5).Value = 5200
$ws.Cells.Item(15, 6).Value = 5100
$ws.Cells.Item(15, 7).Value = 5100
$ws.Cells.Item(15, 8).Value = 5200
$ws.Cells.Item(15, 9).Value = 5200
$ws.Cells.Item(15, 10).Value = 5200
$ws.Cells.Item(15, 11).Value = 5200
$ws.Cells.Item(15, 12).Value = 5200
$ws.Cells.Item(15, 13).Value = 5200
$ws.Cells.Item(17, 4).Value = 2193300
$ws.Cells.Item(17, 5).Value = 1919800
$ws.Cells.Item(17, 6).Value = 1992100
$ws.Cells.Item(17, 7).Value = 1965200
$ws.Cells.Item(17, 8).Value = 1901000
$ws.Cells.Item(17, 9).Value = 1934300
$ws.Cells.Item(17, 10).Value = 1880500
$ws.Cells.Item(17, 11).Value = 1991600
$ws.Cells.Item(17, 12).Value = 1947700
$ws.Cells.Item(17, 13).Value = 1922600
$ws.Cells.Item(18, 4).Value = -263800
$ws.Cells.Item(18, 5).Value = -25700
$ws.Cells.Item(18, 6).Value = -40900
$ws.Cells.Item(18, 7).Value = 15300
$ws.Cells.Item(18, 8).Value = 34000
$ws.Cells.Item(18, 9).Value = 3300
$ws.Cells.Item(18, 10).Value = 46200
$ws.Cells.Item(18, 11).Value = 4100
$ws.Cells.Item(18, 12).Value = 70300
$ws.Cells.Item(18, 13).Value = 42000
$ws.Cells.Item(20, 4).Value = -1200
$ws.Cells.Item(20, 5).Value = -500
$ws.Cells.Item(20, 6).Value = -800
$ws.Cells.Item(20, 7).Value = -500
$ws.Cells.Item(20, 8).Value = -600
$ws.Cells.Item(20, 9).Value = -400
$ws.Cells.Item(20, 10).Value = -200
$ws.Cells.Item(20, 11).Value = -200
$ws.Cells.Item(20, 12).Value = 1400
$ws.Cells.Item(20, 13).Value = 1200
$ws.Cells.Item(21, 4).Value = -227200
$ws.Cells.Item(21, 5).Value = 11900
$ws.Cells.Item(21, 6).Value = -1600
$ws.Cells.Item(21, 7).Value = 55000
$ws.Cells.Item(21, 8).Value = 74400
$ws.Cells.Item(21, 9).Value = 46100
$ws.Cells.Item(21, 10).Value = 89300
$ws.Cells.Item(21, 11).Value = 47100
$ws.Cells.Item(21, 12).Value = 117400
$ws.Cells.Item(21, 13).Value = 88000
$ws.Cells.Item(22, 4).Value = 14500
$ws.Cells.Item(22, 5).Value = 13800
$ws.Cells.Item(22, 6).Value = 14100
$ws.Cells.Item(22, 7).Value = 14000
$ws.Cells.Item(22, 8).Value = 14600
$ws.Cells.Item(22, 9).Value = 16500
$ws.Cells.Item(22, 10).Value = 16400
$ws.Cells.Item(22, 11).Value = 17500
$ws.Cells.Item(22, 12).Value = 16500
$ws.Cells.Item(22, 13).Value = 16600
$ws.Cells.Item(23, 4).Value = -279600
$ws.Cells.Item(23, 5).Value = -40000
$ws.Cells.Item(23, 6).Value = -55700
$ws.Cells.Item(23, 7).Value = 800
$ws.Cells.Item(23, 8).Value = 18900
$ws.Cells.Item(23, 9).Value = -13700
$ws.Cells.Item(23, 10).Value = 29600
$ws.Cells.Item(23, 11).Value = -13600
$ws.Cells.Item(23, 12).Value = 55200
$ws.Cells.Item(23, 13).Value = 26600
$ws.Cells.Item(24, 4).Value = -16300
$ws.Cells.Item(24, 5).Value = -13400
$ws.Cells.Item(24, 6).Value = -13700
$ws.Cells.Item(24, 7).Value = 1100
$ws.Cells.Item(24, 8).Value = 13100
$ws.Cells.Item(24, 9).Value = -3700
$ws.Cells.Item(24, 10).Value = 11900
$ws.Cells.Item(24, 11).Value = -3800
$ws.Cells.Item(24, 12).Value = 21700
$ws.Cells.Item(24, 13).Value = 12100
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 4).Value = -263300
$ws.Cells.Item(26, 5).Value = -26600
$ws.Cells.Item(26, 6).Value = -42000
$ws.Cells.Item(26, 7).Value = -300
$ws.Cells.Item(26, 8).Value = 5800
$ws.Cells.Item(26, 9).Value = -10000
$ws.Cells.Item(26, 10).Value = 17600
$ws.Cells.Item(26, 11).Value = -9800
$ws.Cells.Item(26, 12).Value = 33500
$ws.Cells.Item(26, 13).Value = 14500
$ws.Cells.Item(27, 4).Value = -263100
$ws.Cells.Item(27, 5).Value = -26400
$ws.Cells.Item(27, 6).Value = -42000
$ws.Cells.Item(27, 7).Value = -300
$ws.Cells.Item(27, 8).Value = 5800
$ws.Cells.Item(27, 9).Value = -10000
$ws.Cells.Item(27, 10).Value = 17600
$ws.Cells.Item(27, 11).Value = -9800
$ws.Cells.Item(27, 12).Value = 33500
$ws.Cells.Item(27, 13).Value = 14500
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(29, 4).Value = 3000
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 1900
$ws.Cells.Item(29, 7).Value = "NA"
$ws.Cells.Item(29, 8).Value = 46500
$ws.Cells.Item(29, 9).Value = 11400
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = "NA"
$ws.Cells.Item(29, 12).Value = -700
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(32, 4).Value = 1200
$ws.Cells.Item(32, 5).Value = 500
$ws.Cells.Item(32, 6).Value = 800
$ws.Cells.Item(32, 7).Value = 500
$ws.Cells.Item(32, 8).Value = 600
$ws.Cells.Item(32, 9).Value = 400
$ws.Cells.Item(32, 10).Value = 200
$ws.Cells.Item(32, 11).Value = 200
$ws.Cells.Item(32, 12).Value = -1400
$ws.Cells.Item(32, 13).Value = -1200
$ws.Cells.Item(33, 4).Value = -260100
$ws.Cells.Item(33, 5).Value = -26400
$ws.Cells.Item(33, 6).Value = -40100
$ws.Cells.Item(33, 7).Value = -300
$ws.Cells.Item(33, 8).Value = 52300
$ws.Cells.Item(33, 9).Value = 1400
$ws.Cells.Item(33, 10).Value = 17600
$ws.Cells.Item(33, 11).Value = -9800
$ws.Cells.Item(33, 12).Value = 32800
$ws.Cells.Item(33, 13).Value = 14500
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(35, 4).Value = -260100
$ws.Cells.Item(35, 5).Value = -26400
$ws.Cells.Item(35, 6).Value = -40100
$ws.Cells.Item(35, 7).Value = -300
$ws.Cells.Item(35, 8).Value = 52300
$ws.Cells.Item(35, 9).Value = 1400
$ws.Cells.Item(35, 10).Value = 17600
$ws.Cells.Item(35, 11).Value = -9800
$ws.Cells.Item(35, 12).Value = 32800
$ws.Cells.Item(35, 13).Value = 14500
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(38, 7).Value = 43190
$ws.Cells.Item(38, 8).Value = 43100
$ws.Cells.Item(38, 9).Value = 43008
$ws.Cells.Item(38, 10).Value = 42916
$ws.Cells.Item(38, 11).Value = 42825
$ws.Cells.Item(38, 12).Value = 42735
$ws.Cells.Item(38, 13).Value = 42643
$ws.Cells.Item(41, 4).Value = 24200
$ws.Cells.Item(41, 5).Value = 21800
$ws.Cells.Item(41, 6).Value = 25400
$ws.Cells.Item(41, 7).Value = 28100
$ws.Cells.Item(41, 8).Value = 16500
$ws.Cells.Item(41, 9).Value = 24300
$ws.Cells.Item(41, 10).Value = 31500
$ws.Cells.Item(41, 11).Value = 31600
$ws.Cells.Item(41, 12).Value = 18000
$ws.Cells.Item(41, 13).Value = 28200
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = 0
$ws.Cells.Item(43, 4).Value = 593500
$ws.Cells.Item(43, 5).Value = 616400
$ws.Cells.Item(43, 6).Value = 594700
$ws.Cells.Item(43, 7).Value = 639200
$ws.Cells.Item(43, 8).Value = 678000
$ws.Cells.Item(43, 9).Value = 673800
$ws.Cells.Item(43, 10).Value = 609100
$ws.Cells.Item(43, 11).Value = 651800
$ws.Cells.Item(43, 12).Value = 674800
$ws.Cells.Item(43, 13).Value = 647800
$ws.Cells.Item(44, 4).Value = 255500
$ws.Cells.Item(44, 5).Value = 262200
$ws.Cells.Item(44, 6).Value = 276700
$ws.Cells.Item(44, 7).Value = 287200
$ws.Cells.Item(44, 8).Value = 278100
$ws.Cells.Item(44, 9).Value = 281100
$ws.Cells.Item(44, 10).Value = 287200
$ws.Cells.Item(44, 11).Value = 290900
$ws.Cells.Item(44, 12).Value = 284500
$ws.Cells.Item(44, 13).Value = 261600
$ws.Cells.Item(45, 4).Value = 39100
$ws.Cells.Item(45, 5).Value = 38300
$ws.Cells.Item(45, 6).Value = 46500
$ws.Cells.Item(45, 7).Value = 47500
$ws.Cells.Item(45, 8).Value = 47300
$ws.Cells.Item(45, 9).Value = 37200
$ws.Cells.Item(45, 10).Value = 40300
$ws.Cells.Item(45, 11).Value = 37100
$ws.Cells.Item(45, 12).Value = 81400
$ws.Cells.Item(45, 13).Value = 82200
$ws.Cells.Item(46, 4).Value = 912300
$ws.Cells.Item(46, 5).Value = 938600
$ws.Cells.Item(46, 6).Value = 943300
$ws.Cells.Item(46, 7).Value = 1002100
$ws.Cells.Item(46, 8).Value = 1019900
$ws.Cells.Item(46, 9).Value = 1016400
$ws.Cells.Item(46, 10).Value = 968100
$ws.Cells.Item(46, 11).Value = 1011400
$ws.Cells.Item(46, 12).Value = 1058600
$ws.Cells.Item(46, 13).Value = 1019700
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = 0
$ws.Cells.Item(48, 4).Value = 1006200
$ws.Cells.Item(48, 5).Value = 999400
$ws.Cells.Item(48, 6).Value = 1001800
$ws.Cells.Item(48, 7).Value = 1056000
$ws.Cells.Item(48, 8).Value = 1094100
$ws.Cells.Item(48, 9).Value = 1091400
$ws.Cells.Item(48, 10).Value = 1124100
$ws.Cells.Item(48, 11).Value = 1135700
$ws.Cells.Item(48, 12).Value = 1163900
$ws.Cells.Item(48, 13).Value = 1144700
$ws.Cells.Item(49, 4).Value = 153600
$ws.Cells.Item(49, 5).Value = 349500
$ws.Cells.Item(49, 6).Value = 354600
$ws.Cells.Item(49, 7).Value = 322800
$ws.Cells.Item(49, 8).Value = 328000
$ws.Cells.Item(49, 9).Value = 333200
$ws.Cells.Item(49, 10).Value = 338400
$ws.Cells.Item(49, 11).Value = 320900
$ws.Cells.Item(49, 12).Value = 326000
$ws.Cells.Item(49, 13).Value = 340900
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(52, 4).Value = 46400
$ws.Cells.Item(52, 5).Value = 66900
$ws.Cells.Item(52, 6).Value = 57900
$ws.Cells.Item(52, 7).Value = 63700
$ws.Cells.Item(52, 8).Value = 61900
$ws.Cells.Item(52, 9).Value = 65700
$ws.Cells.Item(52, 10).Value = 63600
$ws.Cells.Item(52, 11).Value = 57800
$ws.Cells.Item(52, 12).Value = 57700
$ws.Cells.Item(52, 13).Value = 62300
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(54, 4).Value = 2118500
$ws.Cells.Item(54, 5).Value = 2354500
$ws.Cells.Item(54, 6).Value = 2357600
$ws.Cells.Item(54, 7).Value = 2444700
$ws.Cells.Item(54, 8).Value = 2503800
$ws.Cells.Item(54, 9).Value = 2506700
$ws.Cells.Item(54, 10).Value = 2494200
$ws.Cells.Item(54, 11).Value = 2525700
$ws.Cells.Item(54, 12).Value = 2606200
$ws.Cells.Item(54, 13).Value = 2567700
$ws.Cells.Item(57, 4).Value = 695300
$ws.Cells.Item(57, 5).Value = 644100
$ws.Cells.Item(57, 6).Value = 641100
$ws.Cells.Item(57, 7).Value = 623800
$ws.Cells.Item(57, 8).Value = 669200
$ws.Cells.Item(57, 9).Value = 647500
$ws.Cells.Item(57, 10).Value = 653300
$ws.Cells.Item(57, 11).Value = 682200
$ws.Cells.Item(57, 12).Value = 458800
$ws.Cells.Item(57, 13).Value = 675900
$ws.Cells.Item(58, 4).Value = 1200
$ws.Cells.Item(58, 5).Value = 1200
$ws.Cells.Item(58, 6).Value = 1200
$ws.Cells.Item(58, 7).Value = 1100
$ws.Cells.Item(58, 8).Value = 1100
$ws.Cells.Item(58, 9).Value = 142900
$ws.Cells.Item(58, 10).Value = 142200
$ws.Cells.Item(58, 11).Value = 141500
$ws.Cells.Item(58, 12).Value = 140800
$ws.Cells.Item(58, 13).Value = 1100
$ws.Cells.Item(59, 4).Value = 4300
$ws.Cells.Item(59, 5).Value = 300
$ws.Cells.Item(59, 6).Value = 1500
$ws.Cells.Item(59, 7).Value = 1400
$ws.Cells.Item(59, 8).Value = 1800
$ws.Cells.Item(59, 9).Value = 800
$ws.Cells.Item(59, 10).Value = 2000
$ws.Cells.Item(59, 11).Value = 2300
$ws.Cells.Item(59, 12).Value = 248200
$ws.Cells.Item(59, 13).Value = 5400
$ws.Cells.Item(60, 4).Value = 700800
$ws.Cells.Item(60, 5).Value = 645500
$ws.Cells.Item(60, 6).Value = 643700
$ws.Cells.Item(60, 7).Value = 626300
$ws.Cells.Item(60, 8).Value = 672200
$ws.Cells.Item(60, 9).Value = 791200
$ws.Cells.Item(60, 10).Value = 797500
$ws.Cells.Item(60, 11).Value = 826000
$ws.Cells.Item(60, 12).Value = 847800
$ws.Cells.Item(60, 13).Value = 682400
$ws.Cells.Item(61, 4).Value = 905200
$ws.Cells.Item(61, 5).Value = 886000
$ws.Cells.Item(61, 6).Value = 855800
$ws.Cells.Item(61, 7).Value = 906100
$ws.Cells.Item(61, 8).Value = 912100
$ws.Cells.Item(61, 9).Value = 803600
$ws.Cells.Item(61, 10).Value = 762100
$ws.Cells.Item(61, 11).Value = 750200
$ws.Cells.Item(61, 12).Value = 745200
$ws.Cells.Item(61, 13).Value = 895000
$ws.Cells.Item(62, 4).Value = 197800
$ws.Cells.Item(62, 5).Value = 238600
$ws.Cells.Item(62, 6).Value = 240600
$ws.Cells.Item(62, 7).Value = 262000
$ws.Cells.Item(62, 8).Value = 263600
$ws.Cells.Item(62, 9).Value = 307500
$ws.Cells.Item(62, 10).Value = 326400
$ws.Cells.Item(62, 11).Value = 354600
$ws.Cells.Item(62, 12).Value = 402600
$ws.Cells.Item(62, 13).Value = 400300
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0
$ws.Cells.Item(66, 4).Value = 1815500
$ws.Cells.Item(66, 5).Value = 1782000
$ws.Cells.Item(66, 6).Value = 1751900
$ws.Cells.Item(66, 7).Value = 1794400
$ws.Cells.Item(66, 8).Value = 1847900
$ws.Cells.Item(66, 9).Value = 1902300
$ws.Cells.Item(66, 10).Value = 1886000
$ws.Cells.Item(66, 11).Value = 1930700
$ws.Cells.Item(66, 12).Value = 1995700
$ws.Cells.Item(66, 13).Value = 1977800
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0
$ws.Cells.Item(72, 4).Value = -261000
$ws.Cells.Item(72, 5).Value = -800
$ws.Cells.Item(72, 6).Value = 33900
$ws.Cells.Item(72, 7).Value = 82400
$ws.Cells.Item(72, 8).Value = 74200
$ws.Cells.Item(72, 9).Value = 30100
$ws.Cells.Item(72, 10).Value = 36900
$ws.Cells.Item(72, 11).Value = 27600
$ws.Cells.Item(72, 12).Value = 45700
$ws.Cells.Item(72, 13).Value = 21100
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0
$ws.Cells.Item(76, 4).Value = 303000
$ws.Cells.Item(76, 5).Value = 572500
$ws.Cells.Item(76, 6).Value = 605700
$ws.Cells.Item(76, 7).Value = 650400
$ws.Cells.Item(76, 8).Value = 655900
$ws.Cells.Item(76, 9).Value = 604400
$ws.Cells.Item(76, 10).Value = 608200
$ws.Cells.Item(76, 11).Value = 595000
$ws.Cells.Item(76, 12).Value = 610600
$ws.Cells.Item(76, 13).Value = 589900
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(80, 7).Value = 43190
$ws.Cells.Item(80, 8).Value = 43100
$ws.Cells.Item(80, 9).Value = 43008
$ws.Cells.Item(80, 10).Value = 42916
$ws.Cells.Item(80, 11).Value = 42825
$ws.Cells.Item(80, 12).Value = 42735
$ws.Cells.Item(80, 13).Value = 42643
$ws.Cells.Item(81, 4).Value = -260100
$ws.Cells.Item(81, 5).Value = -26400
$ws.Cells.Item(81, 6).Value = -40100
$ws.Cells.Item(81, 7).Value = -300
$ws.Cells.Item(81, 8).Value = 52300
$ws.Cells.Item(81, 9).Value = 1400
$ws.Cells.Item(81, 10).Value = 17600
$ws.Cells.Item(81, 11).Value = -9800
$ws.Cells.Item(81, 12).Value = 32800
$ws.Cells.Item(81, 13).Value = 14500
$ws.Cells.Item(83, 4).Value = 37800
$ws.Cells.Item(83, 5).Value = 38100
$ws.Cells.Item(83, 6).Value = 40000
$ws.Cells.Item(83, 7).Value = 40100
$ws.Cells.Item(83, 8).Value = 41000
$ws.Cells.Item(83, 9).Value = 43200
$ws.Cells.Item(83, 10).Value = 43300
$ws.Cells.Item(83, 11).Value = 43200
$ws.Cells.Item(83, 12).Value = 45600
$ws.Cells.Item(83, 13).Value = 44900
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0
$ws.Cells.Item(89, 4).Value = 33200
$ws.Cells.Item(89, 5).Value = -1000
$ws.Cells.Item(89, 6).Value = 81800
$ws.Cells.Item(89, 7).Value = 39000
$ws.Cells.Item(89, 8).Value = 78100
$ws.Cells.Item(89, 9).Value = -12500
$ws.Cells.Item(89, 10).Value = 51600
$ws.Cells.Item(89, 11).Value = 27600
$ws.Cells.Item(89, 12).Value = 72800
$ws.Cells.Item(89, 13).Value = 59300
$ws.Cells.Item(91, 4).Value = -46700
$ws.Cells.Item(91, 5).Value = -31400
$ws.Cells.Item(91, 6).Value = -20800
$ws.Cells.Item(91, 7).Value = -16500
$ws.Cells.Item(91, 8).Value = -45300
$ws.Cells.Item(91, 9).Value = -26800
$ws.Cells.Item(91, 10).Value = -26200
$ws.Cells.Item(91, 11).Value = -8400
$ws.Cells.Item(91, 12).Value = -63300
$ws.Cells.Item(91, 13).Value = -35600
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0
$ws.Cells.Item(94, 4).Value = -46300
$ws.Cells.Item(94, 5).Value = -24700
$ws.Cells.Item(94, 6).Value = -25900
$ws.Cells.Item(94, 7).Value = -12300
$ws.Cells.Item(94, 8).Value = -44100
$ws.Cells.Item(94, 9).Value = -28200
$ws.Cells.Item(94, 10).Value = -55300
$ws.Cells.Item(94, 11).Value = -7400
$ws.Cells.Item(94, 12).Value = -63300
$ws.Cells.Item(94, 13).Value = -32500
$ws.Cells.Item(96, 4).Value = -2700
$ws.Cells.Item(96, 5).Value = -8200
$ws.Cells.Item(96, 6).Value = -8200
$ws.Cells.Item(96, 7).Value = -8200
$ws.Cells.Item(96, 8).Value = -8200
$ws.Cells.Item(96, 9).Value = -8200
$ws.Cells.Item(96, 10).Value = -8200
$ws.Cells.Item(96, 11).Value = -8200
$ws.Cells.Item(96, 12).Value = -8100
$ws.Cells.Item(96, 13).Value = -8200
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0
$ws.Cells.Item(100, 4).Value = 15500
$ws.Cells.Item(100, 5).Value = 22000
$ws.Cells.Item(100, 6).Value = -58600
$ws.Cells.Item(100, 7).Value = -15000
$ws.Cells.Item(100, 8).Value = -41800
$ws.Cells.Item(100, 9).Value = 33500
$ws.Cells.Item(100, 10).Value = 3600
$ws.Cells.Item(100, 11).Value = -6600
$ws.Cells.Item(100, 12).Value = -19000
$ws.Cells.Item(100, 13).Value = -21900
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = -700
$ws.Cells.Item(101, 13).Value = -500
$ws.Cells.Item(102, 4).Value = 2400
$ws.Cells.Item(102, 5).Value = -3600
$ws.Cells.Item(102, 6).Value = -2700
$ws.Cells.Item(102, 7).Value = 11600
$ws.Cells.Item(102, 8).Value = -7800
$ws.Cells.Item(102, 9).Value = -7200
$ws.Cells.Item(102, 10).Value = -100
$ws.Cells.Item(102, 11).Value = 13600
$ws.Cells.Item(102, 12).Value = -10200
$ws.Cells.Item(102, 13).Value = 4300
